$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 111574403
$ws.Range("B14").Value = 89686
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 658
$ws.Range("F14").Value = "Rosenticka"
$ws.Range("G14").Value = "Rhodofomes roseus"
$ws.Range("H14").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q14").Value = 562547.0565141424
$ws.Range("R14").Value = 6954767.535469687
# Row 15
$ws.Range("A15").Value = 111573866
$ws.Range("AB15").Value = "15:17"
$ws.Range("Q15").Value = 562601.7570288588
$ws.Range("R15").Value = 6954814.918206804
$ws.Range("Z15").Value = "15:17"
# Row 16
$ws.Range("A16").Value = 111573533
# Row 17
$ws.Range("A17").Value = 111574338
$ws.Range("AB17").Value = "15:26"
$ws.Range("B17").Value = 89686
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 658
$ws.Range("F17").Value = "Rosenticka"
$ws.Range("G17").Value = "Rhodofomes roseus"
$ws.Range("H17").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 562557.3535548041
$ws.Range("R17").Value = 6954757.635990249
$ws.Range("Z17").Value = "15:26"
# Row 18
$ws.Range("A18").Value = 111573803
$ws.Range("AB18").Value = "15:14"
$ws.Range("B18").Value = 96348
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
$ws.Range("Q18").Value = 562591.0245237258
$ws.Range("R18").Value = 6954847.751526525
$ws.Range("Z18").Value = "15:14"
# Row 19
$ws.Range("A19").Value = 111574509
$ws.Range("AB19").Value = "15:45"
$ws.Range("B19").Value = 96348
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = "Knärot"
$ws.Range("G19").Value = "Goodyera repens"
$ws.Range("H19").Value = "(L.) R. Br."
$ws.Range("M19").ClearContents() | Out-Null
$ws.Range("Q19").Value = 562529.1073683554
$ws.Range("R19").Value = 6954769.030357216
$ws.Range("Z19").Value = "15:45"
# Row 20
$ws.Range("A20").Value = 111574240
$ws.Range("AB20").Value = "15:26"
$ws.Range("AC20").Value = "Familj med 5 talltitor. Permanent revir"
$ws.Range("B20").Value = 56543
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 103021
$ws.Range("F20").Value = "Talltita"
$ws.Range("G20").Value = "Poecile montanus"
$ws.Range("H20").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I20").Value = "'5"
$ws.Range("Q20").Value = 562533.1227179464
$ws.Range("R20").Value = 6954848.029061474
$ws.Range("Z20").Value = "15:26"
# Row 21
$ws.Range("A21").Value = 111573948
$ws.Range("AB21").Value = "15:20"
$ws.Range("Q21").Value = 562576.2301468613
$ws.Range("R21").Value = 6954852.517936011
$ws.Range("Z21").Value = "15:20"
# Row 22
$ws.Range("A22").Value = 111576450
$ws.Range("AB22").Value = "17:10"
$ws.Range("AC22").Value = "Rikligt"
$ws.Range("Q22").Value = 562979.5212303887
$ws.Range("R22").Value = 6954739.97881452
$ws.Range("Z22").Value = "17:10"
# Row 23
$ws.Range("A23").Value = 111574689
$ws.Range("AB23").Value = "15:47"
$ws.Range("AC23").ClearContents() | Out-Null
$ws.Range("Q23").Value = 562517.0252856832
$ws.Range("R23").Value = 6954776.14289257
$ws.Range("Z23").Value = "15:47"
# Row 25
$ws.Range("A25").Value = 111574128
$ws.Range("B25").Value = 96348
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("Q25").Value = 562555.4143375416
$ws.Range("R25").Value = 6954835.60431945
# Row 26
$ws.Range("A26").Value = 111578233
$ws.Range("AB26").Value = "18:43"
$ws.Range("AC26").ClearContents() | Out-Null
$ws.Range("Q26").Value = 563026.0554397166
$ws.Range("R26").Value = 6954541.256262898
$ws.Range("Z26").Value = "18:43"
# Row 27
$ws.Range("A27").Value = 111575868
$ws.Range("AB27").Value = "16:43"
$ws.Range("B27").Value = 96348
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = "Knärot"
$ws.Range("G27").Value = "Goodyera repens"
$ws.Range("H27").Value = "(L.) R. Br."
$ws.Range("Q27").Value = 562854.9195222461
$ws.Range("R27").Value = 6954623.341454657
$ws.Range("Z27").Value = "16:43"
# Row 28
$ws.Range("A28").Value = 111578127
$ws.Range("AB28").Value = "18:30"
$ws.Range("B28").Value = 56543
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 103021
$ws.Range("F28").Value = "Talltita"
$ws.Range("G28").Value = "Poecile montanus"
$ws.Range("H28").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M28").Value = "lockläte, övriga läten"
$ws.Range("Q28").Value = 562937.8205991766
$ws.Range("R28").Value = 6954541.406048392
$ws.Range("Z28").Value = "18:30"
# Row 29
$ws.Range("A29").Value = 111576771
$ws.Range("AB29").Value = "17:24"
$ws.Range("AC29").ClearContents() | Out-Null
$ws.Range("B29").Value = 96348
$ws.Range("D29").Value = "VU"
$ws.Range("E29").Value = 220787
$ws.Range("F29").Value = "Knärot"
$ws.Range("G29").Value = "Goodyera repens"
$ws.Range("H29").Value = "(L.) R. Br."
$ws.Range("I29").ClearContents() | Out-Null
$ws.Range("Q29").Value = 562807.4867926922
$ws.Range("R29").Value = 6954821.585021482
$ws.Range("Z29").Value = "17:24"
# Row 30
$ws.Range("A30").Value = 111576037
$ws.Range("AB30").Value = "16:51"
$ws.Range("B30").Value = 89686
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 658
$ws.Range("F30").Value = "Rosenticka"
$ws.Range("G30").Value = "Rhodofomes roseus"
$ws.Range("H30").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q30").Value = 562852.9463231879
$ws.Range("R30").Value = 6954606.325244571
$ws.Range("Z30").Value = "16:51"
# Row 31
$ws.Range("A31").Value = 111575785
$ws.Range("AB31").Value = "16:39"
$ws.Range("B31").Value = 89845
$ws.Range("D31").Value = "VU"
$ws.Range("E31").Value = 1209
$ws.Range("F31").Value = "Rynkskinn"
$ws.Range("G31").Value = "Phlebia centrifuga"
$ws.Range("H31").Value = "P.Karst."
$ws.Range("Q31").Value = 562859.2727272335
$ws.Range("R31").Value = 6954660.134623887
$ws.Range("Z31").Value = "16:39"
# Row 33
$ws.Range("A33").Value = 111576401
$ws.Range("AB33").Value = "16:51"
$ws.Range("B33").Value = 89369
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 5447
$ws.Range("F33").Value = "Vedticka"
$ws.Range("G33").Value = "Fuscoporia viticola"
$ws.Range("H33").Value = "(Schwein.) Murrill"
$ws.Range("Q33").Value = 562964.914807545
$ws.Range("R33").Value = 6954710.791209211
$ws.Range("Z33").Value = "16:51"
